$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row ranges in column A to their new (shortened) file names.
# The original values were full Windows paths to .pep.xml files; the new
# values are bare .raw file names.
$groups = @(
    @{ Start = 2;   End = 17;  Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H14_100pg_AGC300_1.raw" },
    @{ Start = 19;  End = 34;  Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H14_100pg_AGC300_2.raw" },
    @{ Start = 36;  End = 51;  Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H42_100pg_AGC300_1.raw" },
    @{ Start = 53;  End = 68;  Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H42_100pg_AGC300_2.raw" },
    @{ Start = 70;  End = 85;  Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H98_100pg_AGC300_1.raw" },
    @{ Start = 87;  End = 102; Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H98_100pg_AGC300_2.raw" },
    @{ Start = 104; End = 119; Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H210_100pg_AGC300_1.raw" },
    @{ Start = 121; End = 136; Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H210_100pg_AGC300_2.raw" },
    @{ Start = 138; End = 153; Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H434_100pg_AGC300_1.raw" },
    @{ Start = 155; End = 170; Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_H434_100pg_AGC300_2.raw" },
    @{ Start = 172; End = 187; Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_No126_100pg_AGC300_1.raw" },
    @{ Start = 189; End = 204; Name = "20200820_EXPL8_EVO1_ZY_SA_44min_TMT_No126_100pg_AGC300_2.raw" }
)

foreach ($g in $groups) {
    $rng = $ws.Range("A$($g.Start):A$($g.End)")
    $rng.Value = $g.Name
}

# Shrink column A to fit the new, much shorter file names.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update the active selection to reflect where the author finished editing.
$ws.Range("A6").Select() | Out-Null
